$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 26 (existing rows 26+ shift down to 27+,
# carrying their content/formatting with them - matches target row 27 == old row 26)
$ws.Rows("26:26").Insert()

# --- Update existing rows 24 & 25, and populate the new row 26 ---
# Order below matches the order new shared strings must be appended in.

# Row 25, column F: append the two new related pages to the existing list
$ws.Range("F25").Value = "account_management.php, editor_create_user.php, deactivate_user.php, editor_user_account_management.php"

# Row 24, column F: was a placeholder "NEW?", now points at the real new page
$ws.Range("F24").Value = "editor_find_users.php"

# New row 26: Editor list-of-users page
$ws.Range("B26").Value = "manage\user_accounts.php"
$ws.Range("C26").Value = "Editor list of users"
$ws.Range("D26").Value = "Editor"
$ws.Range("F26").Value = "editor_user_account_management.php, editor_find_users.php"
$ws.Range("G26").Value = "PENDING"

# Row 25, column C: new purpose/notes text
$ws.Range("C25").Value = "Editor making adding or changing a single user *merge pages for processing updates*"

# Match formatting style of neighboring "A" marker cells for the new row
$ws.Range("A26").Style = $ws.Range("A25").Style

# --- Update the sheet view (scroll position / active selection) ---
$ws.Application.ActiveWindow.ScrollRow = 12
$ws.Range("C25").Select()
